$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E7"   = 16.1775
    "B9"   = 6.060599999999997
    "E12"  = 18.11120000000003
    "B18"  = 6.146299999999997
    "B20"  = 9.3004
    "E26"  = 16.1893
    "B27"  = 6.174600000000005
    "E27"  = 16.64069999999998
    "E29"  = 17.11460000000002
    "E37"  = 16.7691
    "E38"  = 16.4389
    "E51"  = 17.20050000000001
    "E55"  = 16.3073
    "B69"  = 5.408099999999995
    "E69"  = 17.30420000000003
    "E70"  = 18.12400000000002
    "B76"  = 5.601099999999997
    "B82"  = 5.158200000000003
    "E83"  = 16.6084
    "E102" = 16.8202
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
